# Commit: "adding new progress as of date 04 nov 2025"
#
# On the "Training Dashboard" sheet, for every data row (3 through 19):
#   - Column H ("PERIOD TO EXPIRE") decreases by 1 day.
#   - Column I ("LAST UPDATE") moves from "03-Nov-2025" to "04-Nov-2025".
#
# Column I stores the date as plain text (General-formatted cell holding an
# inline string), not a real Excel date. Assigning a date-looking string via
# .Value would make Excel auto-convert the cell to a date serial number and
# switch its number format - that is avoided here by writing the text with a
# leading apostrophe (forcing literal text) and then restoring the cell's
# original look-and-feel by copying the format from the untouched neighboring
# cell in column J (which shares the same base style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 19; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H: PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I: LAST UPDATE

    # Decrement the expiry countdown by one day.
    if ($hCell.Value2 -ne $null) {
        $hCell.Value2 = $hCell.Value2 - 1
    }

    # Write the new "last update" date as literal text (apostrophe prefix
    # keeps Excel from reinterpreting it as a date value/format).
    $iCell.Value = "'04-Nov-2025"

    # Restore the plain General formatting/style the cell originally had.
    $formatSource = $ws.Cells.Item($row, 10)  # Column J cell, same base style
    $formatSource.Copy() | Out-Null
    $iCell.PasteSpecial(-4122) | Out-Null     # xlPasteFormats
}

$excel.CutCopyMode = 0
